$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.299.86"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.854.89"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.19"
$ws.Range("E5").Value = "  +2.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7043"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9985"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07731"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3070"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.68"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07817"
$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.42"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.150"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "1.852.18"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6901"
$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.624"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008338"
$ws.Range("E17").Value = "  -1.41%  "

$ws.Range("D18").Value = "29.254.82"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.34"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").Value = "2.095.50"
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.76"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.535"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9976"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.40"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.868"
$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.537"
$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.243"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.195"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05137"
$ws.Range("E33").Value = "  -1.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7917"
$ws.Range("E34").Value = "  +4.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.908"
$ws.Range("E35").Value = "  +3.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.691"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").Value = "1.333.81"
$ws.Range("E38").Value = "  +9.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01875"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.717"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9597"
$ws.Range("E41").Value = "  +6.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.066"
$ws.Range("E42").Value = "  +10.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.20"
$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.725"
$ws.Range("E45").Value = "  +2.06%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "1.994.94"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5185"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.85"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.771"
$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.988"
$ws.Range("E51").Value = "  -0.81%  "
